$d = $word.ActiveDocument

# The commit removes the (redundant / default) <w:contextualSpacing w:val="0"/>
# element from every paragraph's pPr throughout the document body. Word's
# object model has no direct ParagraphFormat.ContextualSpacing property, so
# drive the underlying OOXML directly per-paragraph via WordOpenXML/InsertXML.

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $xml = $p.Range.WordOpenXML
    if ($xml -match '<w:contextualSpacing\b[^>]*/>') {
        $newXml = $xml -replace '<w:contextualSpacing\b[^>]*/>', ''
        $p.Range.InsertXML($newXml)
    }
}

Write-Host "done"
